$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 346021.01
$ws.Range("B3").Value = 340067.57
$ws.Range("B4").Value = 406645.98
$ws.Range("B5").Value = 271806.75
$ws.Range("B6").Value = 426460.18
$ws.Range("B7").Value = 1791001.49
